# Update vm_pu.xlsx results for Case_2_3 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038151829842275
$ws.Range("D2").Value = 1.04505869851694
$ws.Range("E2").Value = 1.053246545436338
$ws.Range("F2").Value = 1.058929423793294
$ws.Range("I2").Value = 1.037420701745724
$ws.Range("J2").Value = 1.043250822113248
$ws.Range("K2").Value = 1.047828004073361
$ws.Range("L2").Value = 1.055993043092899
$ws.Range("M2").Value = 1.06166031593714
$ws.Range("N2").Value = 1.018287057393064
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039253874869203
$ws.Range("D3").Value = 1.045880264076009
$ws.Range("E3").Value = 1.054220975566197
$ws.Range("F3").Value = 1.059900039681179
$ws.Range("I3").Value = 1.037619459843677
$ws.Range("J3").Value = 1.043996776375574
$ws.Range("K3").Value = 1.048460906723066
$ws.Range("L3").Value = 1.056780066166635
$ws.Range("M3").Value = 1.062444665075708
$ws.Range("N3").Value = 1.018539599046028
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039967132775646
$ws.Range("D4").Value = 1.046411895899369
$ws.Range("E4").Value = 1.05485238543627
$ws.Range("F4").Value = 1.060528621524585
$ws.Range("I4").Value = 1.037746758894377
$ws.Range("J4").Value = 1.044479068051359
$ws.Range("K4").Value = 1.04886981203379
$ws.Range("L4").Value = 1.0572895910033
$ws.Range("M4").Value = 1.062952098460688
$ws.Range("N4").Value = 1.01870274988206
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040267025275157
$ws.Range("D5").Value = 1.046635399033182
$ws.Range("E5").Value = 1.055118042106836
$ws.Range("F5").Value = 1.060793003201836
$ws.Range("I5").Value = 1.037799961362054
$ws.Range("J5").Value = 1.044681729758987
$ws.Range("K5").Value = 1.049041565803924
$ws.Range("L5").Value = 1.057503858890709
$ws.Range("M5").Value = 1.063165400683013
$ws.Range("N5").Value = 1.018771275967578
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040317380834286
$ws.Range("D6").Value = 1.046672926474655
$ws.Range("E6").Value = 1.055162659439149
$ws.Range("F6").Value = 1.060837401411175
$ws.Range("I6").Value = 1.037808875866652
$ws.Range("J6").Value = 1.044715752050323
$ws.Range("K6").Value = 1.049070395208961
$ws.Range("L6").Value = 1.057539839127201
$ws.Range("M6").Value = 1.063201213714252
$ws.Range("N6").Value = 1.018782778127614
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039971139800745
$ws.Range("D7").Value = 1.046414882339866
$ws.Range("E7").Value = 1.054855934321752
$ws.Range("F7").Value = 1.060532153711719
$ws.Range("I7").Value = 1.037747471022694
$ws.Range("J7").Value = 1.044481776395492
$ws.Range("K7").Value = 1.048872107606049
$ws.Range("L7").Value = 1.057292453811941
$ws.Range("M7").Value = 1.062954948705555
$ws.Range("N7").Value = 1.018703665776855
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038524237952697
$ws.Range("D8").Value = 1.045336345106742
$ws.Range("E8").Value = 1.053575674270584
$ws.Range("F8").Value = 1.059257338168413
$ws.Range("I8").Value = 1.037488144213014
$ws.Range("J8").Value = 1.043503001607683
$ws.Range("K8").Value = 1.048042025668014
$ws.Range("L8").Value = 1.05625896532678
$ws.Range("M8").Value = 1.061925409335924
$ws.Range("N8").Value = 1.018372458868414
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035975831420586
$ws.Range("D9").Value = 1.043436035458829
$ws.Range("E9").Value = 1.0513265285944
$ws.Range("F9").Value = 1.057015025818063
$ws.Range("I9").Value = 1.037021150075941
$ws.Range("J9").Value = 1.041775292196236
$ws.Range("K9").Value = 1.046574545226015
$ws.Range("L9").Value = 1.054439906697188
$ws.Range("M9").Value = 1.060110540121086
$ws.Range("N9").Value = 1.017786841851843
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034277690733791
$ws.Range("D10").Value = 1.042169338762099
$ws.Range("E10").Value = 1.049831735530615
$ws.Range("F10").Value = 1.055522926925956
$ws.Range("I10").Value = 1.036703093587203
$ws.Range("J10").Value = 1.040621489343516
$ws.Range("K10").Value = 1.045593038865726
$ws.Range("L10").Value = 1.053228632099609
$ws.Range("M10").Value = 1.058900192476281
$ws.Range("N10").Value = 1.017395102484092
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033542560121018
$ws.Range("D11").Value = 1.041620893635873
$ws.Range("E11").Value = 1.049185581654331
$ws.Range("F11").Value = 1.054877496954948
$ws.Range("I11").Value = 1.036563778852135
$ws.Range("J11").Value = 1.040121407599226
$ws.Range("K11").Value = 1.04516728391952
$ws.Range("L11").Value = 1.05270448260906
$ws.Range("M11").Value = 1.058376000862988
$ws.Range("N11").Value = 1.017225161607419
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033269525623616
$ws.Range("D12").Value = 1.04141718360493
$ws.Range("E12").Value = 1.048945737394036
$ws.Range("F12").Value = 1.054637854965315
$ws.Range("I12").Value = 1.03651179182651
$ws.Range("J12").Value = 1.039935583236568
$ws.Range("K12").Value = 1.045009026192575
$ws.Range("L12").Value = 1.052509841738105
$ws.Range("M12").Value = 1.058181277628473
$ws.Range("N12").Value = 1.017161990747365
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033328091295177
$ws.Range("D13").Value = 1.041460879778036
$ws.Range("E13").Value = 1.048997177278049
$ws.Range("F13").Value = 1.054689254478063
$ws.Range("I13").Value = 1.036522954047318
$ws.Range("J13").Value = 1.039975446439357
$ws.Range("K13").Value = 1.045042978139364
$ws.Range("L13").Value = 1.052551590524041
$ws.Range("M13").Value = 1.058223047102159
$ws.Range("N13").Value = 1.017175543248007
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03351999048272
$ws.Range("D14").Value = 1.041604054757563
$ws.Range("E14").Value = 1.049165752661988
$ws.Range("F14").Value = 1.05485768604203
$ws.Range("I14").Value = 1.036559486469612
$ws.Range("J14").Value = 1.040106048762833
$ws.Range("K14").Value = 1.045154204598095
$ws.Range("L14").Value = 1.052688392478996
$ws.Range("M14").Value = 1.05835990528624
$ws.Range("N14").Value = 1.017219940847794
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033638229298335
$ws.Range("D15").Value = 1.041692270517088
$ws.Range("E15").Value = 1.049269639541819
$ws.Range("F15").Value = 1.054961475478514
$ws.Range("I15").Value = 1.03658196359247
$ws.Range("J15").Value = 1.040186507659007
$ws.Range("K15").Value = 1.045222719872757
$ws.Range("L15").Value = 1.05277268752768
$ws.Range("M15").Value = 1.058444226136267
$ws.Range("N15").Value = 1.017247289415844
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034326482575875
$ws.Range("D16").Value = 1.042205738202302
$ws.Range("E16").Value = 1.049874641890483
$ws.Range("F16").Value = 1.055565775906197
$ws.Range("I16").Value = 1.036712305874962
$ws.Range("J16").Value = 1.040654668061741
$ws.Range("K16").Value = 1.045621278906322
$ws.Range("L16").Value = 1.053263425385007
$ws.Range("M16").Value = 1.058934979205816
$ws.Range("N16").Value = 1.017406374268773
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034758252220894
$ws.Range("D17").Value = 1.042527834729034
$ws.Range("E17").Value = 1.050254439167447
$ws.Range("F17").Value = 1.055945014288006
$ws.Range("I17").Value = 1.036793639463517
$ws.Range("J17").Value = 1.040948204875733
$ws.Range("K17").Value = 1.045871082160179
$ws.Range("L17").Value = 1.053571343728187
$ws.Range("M17").Value = 1.059242788292768
$ws.Range("N17").Value = 1.01750607968806
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035010113125821
$ws.Range("D18").Value = 1.042715712292123
$ws.Range("E18").Value = 1.050476074721585
$ws.Range("F18").Value = 1.056166281142451
$ws.Range("I18").Value = 1.036840926142818
$ws.Range("J18").Value = 1.041119373861583
$ws.Range("K18").Value = 1.046016715192607
$ws.Range("L18").Value = 1.053750980102283
$ws.Range("M18").Value = 1.05942231815409
$ws.Range("N18").Value = 1.017564205758935
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035095994084197
$ws.Range("D19").Value = 1.042779774368717
$ws.Range("E19").Value = 1.050551664718193
$ws.Range("F19").Value = 1.056241738181756
$ws.Range("I19").Value = 1.036857023581392
$ws.Range("J19").Value = 1.041177730247794
$ws.Range("K19").Value = 1.046066359907649
$ws.Range("L19").Value = 1.053812237015907
$ws.Range("M19").Value = 1.059483531495088
$ws.Range("N19").Value = 1.017584020085463
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034711925672667
$ws.Range("D20").Value = 1.042493276391621
$ws.Range("E20").Value = 1.050213679521994
$ws.Range("F20").Value = 1.055904319024743
$ws.Range("I20").Value = 1.036784929048919
$ws.Range("J20").Value = 1.040916715918879
$ws.Range("K20").Value = 1.045844288181999
$ws.Range("L20").Value = 1.053538303606671
$ws.Range("M20").Value = 1.059209764319526
$ws.Range("N20").Value = 1.017495385385245
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033463480249292
$ws.Range("D21").Value = 1.041561893100216
$ws.Range("E21").Value = 1.049116106824268
$ws.Range("F21").Value = 1.054808084395637
$ws.Range("I21").Value = 1.036548735187047
$ws.Range("J21").Value = 1.040067591614013
$ws.Range("K21").Value = 1.045121454306331
$ws.Range("L21").Value = 1.052648106285309
$ws.Range("M21").Value = 1.058319604376895
$ws.Range("N21").Value = 1.017206868164857
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03267868170056
$ws.Range("D22").Value = 1.040976335657027
$ws.Range("E22").Value = 1.048426980366647
$ws.Range("F22").Value = 1.054119413737135
$ws.Range("I22").Value = 1.0363988461302
$ws.Range("J22").Value = 1.039533298417084
$ws.Range("K22").Value = 1.044666324073229
$ws.Range("L22").Value = 1.052088702588455
$ws.Range("M22").Value = 1.057759838133026
$ws.Range("N22").Value = 1.017025192599501
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033094704378513
$ws.Range("D23").Value = 1.041286746772172
$ws.Range("E23").Value = 1.048792208046498
$ws.Range("F23").Value = 1.054484436367079
$ws.Range("I23").Value = 1.036478436340096
$ws.Range("J23").Value = 1.039816576690666
$ws.Range("K23").Value = 1.044907659314904
$ws.Range("L23").Value = 1.05238522461104
$ws.Range("M23").Value = 1.058056588970471
$ws.Range("N23").Value = 1.017121528131688
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034732858592994
$ws.Range("D24").Value = 1.042508891804016
$ws.Range("E24").Value = 1.050232096723053
$ws.Range("F24").Value = 1.055922707266187
$ws.Range("I24").Value = 1.036788865385615
$ws.Range("J24").Value = 1.04093094456529
$ws.Range("K24").Value = 1.045856395452829
$ws.Range("L24").Value = 1.053553232913969
$ws.Range("M24").Value = 1.059224686462289
$ws.Range("N24").Value = 1.017500217774303
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036634513453806
$ws.Range("D25").Value = 1.043927282495904
$ws.Range("E25").Value = 1.051907172306743
$ws.Range("F25").Value = 1.05759423064282
$ws.Range("I25").Value = 1.037143066019343
$ws.Range("J25").Value = 1.04222229907313
$ws.Range("K25").Value = 1.046954487199109
$ws.Range("L25").Value = 1.054909927091071
$ws.Range("M25").Value = 1.060579806052637
$ws.Range("N25").Value = 1.017938472594447

Write-Host "Updated $(264) cells in vm_pu sheet"
